$d = $word.ActiveDocument

# Rebuilds a paragraph's runs from scratch (used to split a single run's
# text into several runs with independent text, as a plain Find/Replace
# cannot introduce new <w:r> boundaries). The paragraph's own attributes
# (w14:paraId, rsids, ...) are read back off the live paragraph so they
# survive the round trip unchanged.
function Set-ParagraphRuns {
    param(
        [int]$ParaIndex,
        [string[]]$Texts
    )

    $p = $d.Paragraphs($ParaIndex)
    $r = $p.Range

    $paraAttrs = ""
    $existingXml = $r.WordOpenXML
    $m = [regex]::Match($existingXml, '<w:p\s+([^>]*?)/?>')
    if ($m.Success) {
        $paraAttrs = $m.Groups[1].Value
    }

    $runsXml = ""
    foreach ($t in $Texts) {
        $needsPreserve = ($t.Length -eq 0) -or ($t.StartsWith(" ")) -or ($t.EndsWith(" "))
        if ($needsPreserve) {
            $runsXml += '<w:r><w:t xml:space="preserve">' + $t + '</w:t></w:r>'
        } else {
            $runsXml += '<w:r><w:t>' + $t + '</w:t></w:r>'
        }
    }

    $xml = '<?xml version="1.0" standalone="yes"?>' +
           '<?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
             '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData>' +
                 '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
                   '<w:body>' +
                     '<w:p ' + $paraAttrs + '>' + $runsXml + '</w:p>' +
                   '</w:body>' +
                 '</w:document>' +
               '</pkg:xmlData>' +
             '</pkg:part>' +
           '</pkg:package>'

    # When targeting the very last paragraph in the story, its Range.End
    # coincides with the end of the document content; replacing that
    # exact span via InsertXML leaves a stray trailing empty paragraph
    # behind. Shrink the span by one char (the paragraph mark) to dodge it.
    $docEnd = $d.Content.End
    if ($r.End -eq $docEnd) {
        $r = $d.Range($r.Start, $r.End - 1)
    }

    $r.InsertXML($xml)
}

Set-ParagraphRuns 1 @("Nombre Completo ", "- {{ nombre_completo }}")
Set-ParagraphRuns 2 @("# de Cedula ", "- {{ cedula_ciudadania }}")
Set-ParagraphRuns 3 @("RH ", " - {{ tipo_sangre }}")
Set-ParagraphRuns 4 @("Cargo", " - {{ cargo }}")
Set-ParagraphRuns 5 @("Numero Telefonico", " -", " ", "{{ telefono }}")
Set-ParagraphRuns 6 @("Direccion ", "– {{ dirección_residencia }}")
Set-ParagraphRuns 7 @("Correo Electronico ", "– {{ correo_electronico }}")
Set-ParagraphRuns 8 @("Fecha Nacimiento ", "– {{ fecha_nacimiento }}")

Write-Output "All paragraphs updated."
